# Apply the "shift down + new catch row" edit to rows 18-23 on the active sheet.
#
# Semantics (derived from the target diff):
#   - Row 18 becomes a new "catch" row (stimuli/catch_18.jpg) and loses its
#     category/target/score columns (H, I, M:V).
#   - Rows 19-23 each take on the H:V content that used to belong to the row
#     directly above them (rows 18-22 respectively).
#   - Row 24 onward is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current H:V values for rows 18-22 (old data), since we need the
# pre-edit values of row N to write into row N+1.
$cols = @("H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$snapshot = @{}
for ($r = 18; $r -le 22; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# --- Row 18: clear H, I, M:V; set J/K/L to the new catch values ---
$ws.Range("H18").ClearContents()
$ws.Range("I18").ClearContents()
$ws.Range("J18").Value = "catch"
$ws.Range("K18").Value = "f"
$ws.Range("L18").Value = "stimuli/catch_18.jpg"
$ws.Range("M18:V18").ClearContents()

# --- Rows 19-23: take H:V content from the row above (pre-edit snapshot) ---
for ($r = 19; $r -le 23; $r++) {
    $src = $snapshot[$r - 1]
    foreach ($c in $cols) {
        $val = $src[$c]
        if ($null -eq $val -or $val -eq "") {
            $ws.Range("$c$r").ClearContents()
        } else {
            $ws.Range("$c$r").Value = $val
        }
    }
}
